$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "walkingToRunning"
$ws.Cells.Item(2, 3).Value = 13.38702964782715
$ws.Cells.Item(2, 4).Value = -17.40622901916504
$ws.Cells.Item(2, 5).Value = -7.910804748535156
$ws.Cells.Item(2, 6).Value = -1.239354553595497
$ws.Cells.Item(2, 7).Value = 0.5934752043385932
$ws.Cells.Item(2, 8).Value = 1.315698351103824

$ws.Cells.Item(3, 1).Value = 100
$ws.Cells.Item(3, 2).Value = "walkingToRunning"
$ws.Cells.Item(3, 3).Value = 5.027063369750977
$ws.Cells.Item(3, 4).Value = -9.992935180664062
$ws.Cells.Item(3, 5).Value = -3.668930053710938
$ws.Cells.Item(3, 6).Value = 0.1265866749218411
$ws.Cells.Item(3, 7).Value = 0.9830147470409047
$ws.Cells.Item(3, 8).Value = 1.104435794153933

$ws.Cells.Item(4, 1).Value = 200
$ws.Cells.Item(4, 2).Value = "walkingToRunning"
$ws.Cells.Item(4, 3).Value = 3.064489126205444
$ws.Cells.Item(4, 4).Value = -10.04933261871338
$ws.Cells.Item(4, 5).Value = -9.655683517456056
$ws.Cells.Item(4, 6).Value = 0.2933733014283093
$ws.Cells.Item(4, 7).Value = 0.4318224985956808
$ws.Cells.Item(4, 8).Value = 1.679565984772173

$ws.Cells.Item(5, 1).Value = 300
$ws.Cells.Item(5, 2).Value = "walkingToRunning"
$ws.Cells.Item(5, 3).Value = -11.23465538024902
$ws.Cells.Item(5, 4).Value = -9.59632682800293
$ws.Cells.Item(5, 5).Value = -12.52409744262695
$ws.Cells.Item(5, 6).Value = -0.3514636487162646
$ws.Cells.Item(5, 7).Value = 1.303912767754777
$ws.Cells.Item(5, 8).Value = 1.717408596215416

$ws.Cells.Item(6, 1).Value = 400
$ws.Cells.Item(6, 2).Value = "walkingToRunning"
$ws.Cells.Item(6, 3).Value = 0.3307132720947265
$ws.Cells.Item(6, 4).Value = -10.98133277893066
$ws.Cells.Item(6, 5).Value = -7.000489711761475
$ws.Cells.Item(6, 6).Value = -0.8439442796328934
$ws.Cells.Item(6, 7).Value = 1.107491066802433
$ws.Cells.Item(6, 8).Value = -0.297189421712284

$ws.Cells.Item(7, 1).Value = 500
$ws.Cells.Item(7, 2).Value = "walkingToRunning"
$ws.Cells.Item(7, 3).Value = 15.33979511260986
$ws.Cells.Item(7, 4).Value = -22.00781440734864
$ws.Cells.Item(7, 5).Value = -10.50343608856201
$ws.Cells.Item(7, 6).Value = 1.054603315660084
$ws.Cells.Item(7, 7).Value = -0.3611965101983832
$ws.Cells.Item(7, 8).Value = -2.38295376773451

$ws.Cells.Item(8, 1).Value = 600
$ws.Cells.Item(8, 2).Value = "walkingToRunning"
$ws.Cells.Item(8, 3).Value = 2.999754428863525
$ws.Cells.Item(8, 4).Value = -15.48295783996582
$ws.Cells.Item(8, 5).Value = 6.745935916900635
$ws.Cells.Item(8, 6).Value = 1.400082142342559
$ws.Cells.Item(8, 7).Value = -1.373465267046946
$ws.Cells.Item(8, 8).Value = -4.113783430948127

$ws.Cells.Item(9, 1).Value = 700
$ws.Cells.Item(9, 2).Value = "walkingToRunning"
$ws.Cells.Item(9, 3).Value = -6.992471218109131
$ws.Cells.Item(9, 4).Value = 1.269168853759766
$ws.Cells.Item(9, 5).Value = -5.262283325195312
$ws.Cells.Item(9, 6).Value = 2.811405012786107
$ws.Cells.Item(9, 7).Value = 2.239890750284216
$ws.Cells.Item(9, 8).Value = -2.48325762769725

$ws.Cells.Item(10, 1).Value = 800
$ws.Cells.Item(10, 2).Value = "walkingToRunning"
$ws.Cells.Item(10, 3).Value = 1.578595161437989
$ws.Cells.Item(10, 4).Value = -59.7728385925293
$ws.Cells.Item(10, 5).Value = 11.89765167236328
$ws.Cells.Item(10, 6).Value = 1.144158274066543
$ws.Cells.Item(10, 7).Value = 0.4079549139291005
$ws.Cells.Item(10, 8).Value = -2.216124885407739

$ws.Cells.Item(11, 1).Value = 900
$ws.Cells.Item(11, 2).Value = "walkingToRunning"
$ws.Cells.Item(11, 3).Value = -58.78142929077149
$ws.Cells.Item(11, 4).Value = -6.325687885284424
$ws.Cells.Item(11, 5).Value = -18.43083572387696
$ws.Cells.Item(11, 6).Value = -3.973825000981487
$ws.Cells.Item(11, 7).Value = -4.923421618171759
$ws.Cells.Item(11, 8).Value = 5.359984481912402

$ws.Cells.Item(12, 1).Value = 1000
$ws.Cells.Item(12, 2).Value = "walkingToRunning"
$ws.Cells.Item(12, 3).Value = 1.685664176940918
$ws.Cells.Item(12, 4).Value = -10.93907737731934
$ws.Cells.Item(12, 5).Value = -4.804312705993652
$ws.Cells.Item(12, 6).Value = -10.63765511113626
$ws.Cells.Item(12, 7).Value = -7.195822839694799
$ws.Cells.Item(12, 8).Value = 9.095965221589761

$ws.Cells.Item(13, 1).Value = 1100
$ws.Cells.Item(13, 2).Value = "walkingToRunning"
$ws.Cells.Item(13, 3).Value = 20.54140090942383
$ws.Cells.Item(13, 4).Value = -30.76116371154785
$ws.Cells.Item(13, 5).Value = -21.73070907592773
$ws.Cells.Item(13, 6).Value = 3.433085498305679
$ws.Cells.Item(13, 7).Value = 4.357403949493808
$ws.Cells.Item(13, 8).Value = -2.281185732013872

$ws.Cells.Item(14, 1).Value = 1200
$ws.Cells.Item(14, 2).Value = "walkingToRunning"
$ws.Cells.Item(14, 3).Value = -8.643145561218262
$ws.Cells.Item(14, 4).Value = -12.25782108306885
$ws.Cells.Item(14, 5).Value = 13.20825862884522
$ws.Cells.Item(14, 6).Value = 4.477401155732298
$ws.Cells.Item(14, 7).Value = -3.086138354523986
$ws.Cells.Item(14, 8).Value = -4.067132330150736

$ws.Cells.Item(15, 1).Value = 1300
$ws.Cells.Item(15, 2).Value = "walkingToRunning"
$ws.Cells.Item(15, 3).Value = -1.923628926277161
$ws.Cells.Item(15, 4).Value = -1.592215418815613
$ws.Cells.Item(15, 5).Value = -5.497756004333496
$ws.Cells.Item(15, 6).Value = 5.366742232297456
$ws.Cells.Item(15, 7).Value = 3.331075757610629
$ws.Cells.Item(15, 8).Value = 1.803214041672062

$ws.Cells.Item(16, 1).Value = 1400
$ws.Cells.Item(16, 2).Value = "walkingToRunning"
$ws.Cells.Item(16, 3).Value = 23.09793281555176
$ws.Cells.Item(16, 4).Value = -65.25662231445312
$ws.Cells.Item(16, 5).Value = 1.447998523712158
$ws.Cells.Item(16, 6).Value = 0.746474412569362
$ws.Cells.Item(16, 7).Value = 1.614592228954578
$ws.Cells.Item(16, 8).Value = -2.489844009214525

$ws.Cells.Item(17, 1).Value = 1500
$ws.Cells.Item(17, 2).Value = "walkingToRunning"
$ws.Cells.Item(17, 3).Value = 13.40246772766113
$ws.Cells.Item(17, 4).Value = -2.679043769836426
$ws.Cells.Item(17, 5).Value = -17.18461418151855
$ws.Cells.Item(17, 6).Value = -4.899446143452909
$ws.Cells.Item(17, 7).Value = 0.4346358306082141
$ws.Cells.Item(17, 8).Value = 4.76701935377427

$ws.Cells.Item(18, 1).Value = 1600
$ws.Cells.Item(18, 2).Value = "walkingToRunning"
$ws.Cells.Item(18, 3).Value = -13.39319038391113
$ws.Cells.Item(18, 4).Value = -15.08681869506836
$ws.Cells.Item(18, 5).Value = -24.70808029174805
$ws.Cells.Item(18, 6).Value = -9.7976465981438
$ws.Cells.Item(18, 7).Value = 3.087984385994603
$ws.Cells.Item(18, 8).Value = -1.579483953341611

$ws.Cells.Item(19, 1).Value = 1700
$ws.Cells.Item(19, 2).Value = "walkingToRunning"
$ws.Cells.Item(19, 3).Value = 47.38216018676758
$ws.Cells.Item(19, 4).Value = -81.48814392089844
$ws.Cells.Item(19, 5).Value = 17.70297622680664
$ws.Cells.Item(19, 6).Value = -0.5613160348673585
$ws.Cells.Item(19, 7).Value = 0.9123670508682955
$ws.Cells.Item(19, 8).Value = 0.6450322871691367

$ws.Cells.Item(20, 1).Value = 1800
$ws.Cells.Item(20, 2).Value = "walkingToRunning"
$ws.Cells.Item(20, 3).Value = -2.278211832046509
$ws.Cells.Item(20, 4).Value = 0.415550947189331
$ws.Cells.Item(20, 5).Value = -10.76592063903809
$ws.Cells.Item(20, 6).Value = 2.186993496270989
$ws.Cells.Item(20, 7).Value = -7.507715720437281
$ws.Cells.Item(20, 8).Value = 1.97455261667396

$ws.Cells.Item(21, 1).Value = 1900
$ws.Cells.Item(21, 2).Value = "walkingToRunning"
$ws.Cells.Item(21, 3).Value = 13.07781982421875
$ws.Cells.Item(21, 4).Value = -39.51717376708984
$ws.Cells.Item(21, 5).Value = 42.08852005004883
$ws.Cells.Item(21, 6).Value = 2.968808529135423
$ws.Cells.Item(21, 7).Value = 8.113890336998949
$ws.Cells.Item(21, 8).Value = -6.538449691780965

$ws.Cells.Item(22, 1).Value = 2000
$ws.Cells.Item(22, 2).Value = "walkingToRunning"
$ws.Cells.Item(22, 3).Value = -72.39402770996094
$ws.Cells.Item(22, 4).Value = -27.50520896911621
$ws.Cells.Item(22, 5).Value = -27.84894561767578
$ws.Cells.Item(22, 6).Value = -1.87808840600416
$ws.Cells.Item(22, 7).Value = -1.777517291417992
$ws.Cells.Item(22, 8).Value = 1.215824003261815

$ws.Cells.Item(23, 1).Value = 2100
$ws.Cells.Item(23, 2).Value = "walkingToRunning"
$ws.Cells.Item(23, 3).Value = -8.723164558410645
$ws.Cells.Item(23, 4).Value = -8.500687599182129
$ws.Cells.Item(23, 5).Value = -15.66308975219727
$ws.Cells.Item(23, 6).Value = -4.240441676278492
$ws.Cells.Item(23, 7).Value = 2.87217904187515
$ws.Cells.Item(23, 8).Value = 8.334868672660862

$ws.Cells.Item(24, 1).Value = 2200
$ws.Cells.Item(24, 2).Value = "walkingToRunning"
$ws.Cells.Item(24, 3).Value = -7.880284786224365
$ws.Cells.Item(24, 4).Value = -21.83874320983887
$ws.Cells.Item(24, 5).Value = -8.497885704040527
$ws.Cells.Item(24, 6).Value = 1.2107466585311
$ws.Cells.Item(24, 7).Value = 12.91206166397641
$ws.Cells.Item(24, 8).Value = -4.052727646764785

$ws.Cells.Item(25, 1).Value = 2300
$ws.Cells.Item(25, 2).Value = "walkingToRunning"
$ws.Cells.Item(25, 3).Value = -9.739827156066896
$ws.Cells.Item(25, 4).Value = -0.5102891325950623
$ws.Cells.Item(25, 5).Value = -2.436953544616699
$ws.Cells.Item(25, 6).Value = -0.8707238068140057
$ws.Cells.Item(25, 7).Value = -7.102825836988711
$ws.Cells.Item(25, 8).Value = 0.5642036477901016

$ws.Cells.Item(26, 1).Value = 2400
$ws.Cells.Item(26, 2).Value = "walkingToRunning"
$ws.Cells.Item(26, 3).Value = -2.781617164611816
$ws.Cells.Item(26, 4).Value = 10.22575092315674
$ws.Cells.Item(26, 5).Value = -2.437583446502685
$ws.Cells.Item(26, 6).Value = 6.469962025529018
$ws.Cells.Item(26, 7).Value = -3.154151007992533
$ws.Cells.Item(26, 8).Value = -3.828252468865385

$ws.Cells.Item(27, 1).Value = 2500
$ws.Cells.Item(27, 2).Value = "walkingToRunning"
$ws.Cells.Item(27, 3).Value = 34.05283355712891
$ws.Cells.Item(27, 4).Value = -58.47340393066406
$ws.Cells.Item(27, 5).Value = 20.97058296203613
$ws.Cells.Item(27, 6).Value = -1.114345516927603
$ws.Cells.Item(27, 7).Value = -1.64307104963562
$ws.Cells.Item(27, 8).Value = 0.5952699299951911

$ws.Cells.Item(28, 1).Value = 2600
$ws.Cells.Item(28, 2).Value = "walkingToRunning"
$ws.Cells.Item(28, 3).Value = -11.52477169036865
$ws.Cells.Item(28, 4).Value = 17.23179626464844
$ws.Cells.Item(28, 5).Value = -1.970066547393799
$ws.Cells.Item(28, 6).Value = -2.965141663992466
$ws.Cells.Item(28, 7).Value = -1.430339549081439
$ws.Cells.Item(28, 8).Value = 1.28920519351959

$ws.Cells.Item(29, 1).Value = 2700
$ws.Cells.Item(29, 2).Value = "walkingToRunning"
$ws.Cells.Item(29, 3).Value = -24.07291412353516
$ws.Cells.Item(29, 4).Value = -20.21802139282227
$ws.Cells.Item(29, 5).Value = -17.61478233337402
$ws.Cells.Item(29, 6).Value = -4.759042248326804
$ws.Cells.Item(29, 7).Value = -5.845142855518262
$ws.Cells.Item(29, 8).Value = -5.263579006762336

$ws.Cells.Item(30, 1).Value = 2800
$ws.Cells.Item(30, 2).Value = "walkingToRunning"
$ws.Cells.Item(30, 3).Value = 18.03611946105957
$ws.Cells.Item(30, 4).Value = -66.51955413818359
$ws.Cells.Item(30, 5).Value = 48.28945922851562
$ws.Cells.Item(30, 6).Value = 3.102689377537001
$ws.Cells.Item(30, 7).Value = -4.462434116439107
$ws.Cells.Item(30, 8).Value = 5.959941364070088

$ws.Cells.Item(31, 1).Value = 2900
$ws.Cells.Item(31, 2).Value = "walkingToRunning"
$ws.Cells.Item(31, 3).Value = 12.77012634277344
$ws.Cells.Item(31, 4).Value = 6.439512729644775
$ws.Cells.Item(31, 5).Value = -10.52412605285644
$ws.Cells.Item(31, 6).Value = 3.054766318871562
$ws.Cells.Item(31, 7).Value = -4.391178798045361
$ws.Cells.Item(31, 8).Value = 1.19998335943338

